$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1) Status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet: zh-cn / de-de status cells (row 2, columns E/F)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn / de-de detail sheets: Status column (column C, row 2)
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2) Narrow the Status-related columns ---
# Target stored width ~13.41 "characters" (was ~17.22). ColumnWidth = 12.5 is the
# closest value this engine's pixel-rounded ColumnWidth can reach.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
